$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "Sirve para contar filas, columnas.",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Sirve para ver la cantidad de registros que hay en una fila.",
    2
)
